# ForecastReader: pulls the pilot-turnover forecast numbers and writes
# them into a new "Pilotaż obrotu" column (D) on the report sheet.
# NOTE: bugged — the reader hands back its values as text, so they land
# in the sheet as strings instead of numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$forecast = @(
    "126669", "124994", "99470",  "4770",   "137076", "130174", "264911", "888064",
    "7125",   "11106",  "140474", "142267", "133568", "153107", "256672", "844319",
    "7840",   "136106", "121366", "113134", "121171", "115064", "169608", "784289",
    "10002",  "85496",  "73450",  "74707",  "68280",  "76074",  "115731"
)

# Header for the new column.
$ws.Range("D2").Value = "Pilotaż obrotu"

# Match column D's width to the rest of the bestFit columns.
$ws.Columns.Item(4).ColumnWidth = 12.75

# Write the forecast series into D4:D34 (one value per report day row),
# keeping column B's number format on the cells. The ForecastReader bug:
# values are pushed in as text (quoted) rather than as real numbers, so
# we force a text format before the write, then restore the display
# format afterwards — the cell keeps storing a string.
$row = 4
foreach ($value in $forecast) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "#.###"
    $row = $row + 1
}
